$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "89.350.51"
Set-TextValue $ws.Range("E2") "  -1.12%  "

# Row 3
Set-TextValue $ws.Range("D3") "3.133.34"
Set-TextValue $ws.Range("E3") "  -1.23%  "

# Row 4
Set-TextValue $ws.Range("E4") "  +0.07%  "

# Row 5
Set-TextValue $ws.Range("D5") "214.66"
Set-TextValue $ws.Range("E5") "  +1.22%  "

# Row 6
Set-TextValue $ws.Range("D6") "635.47"
Set-TextValue $ws.Range("E6") "  +3.73%  "

# Row 7
Set-TextValue $ws.Range("D7") "0.393"
Set-TextValue $ws.Range("E7") "  +1.62%  "

# Row 8
Set-TextValue $ws.Range("D8") "0.795"
Set-TextValue $ws.Range("E8") "  +16.14%  "

# Row 9
Set-TextValue $ws.Range("E9") "  +0.08%  "

# Row 10
Set-TextValue $ws.Range("D10") "3.131.00"
Set-TextValue $ws.Range("E10") "  -1.54%  "

# Row 11
Set-TextValue $ws.Range("D11") "0.562"
Set-TextValue $ws.Range("E11") "  -0.85%  "

# Row 12
Set-TextValue $ws.Range("E12") "  +1.71%  "

# Row 13
Set-TextValue $ws.Range("D13") "0.0000251"
Set-TextValue $ws.Range("E13") "  -1.01%  "

# Row 14
Set-TextValue $ws.Range("E14") "  +2.91%  "

# Row 15
Set-TextValue $ws.Range("D15") "89.165.37"
Set-TextValue $ws.Range("E15") "  -1.19%  "

# Row 16
Set-TextValue $ws.Range("D16") "3.707.14"
Set-TextValue $ws.Range("E16") "  -4.02%  "

# Row 17
Set-TextValue $ws.Range("D17") "32.19"
Set-TextValue $ws.Range("E17") "  -1.08%  "

# Row 18
Set-TextValue $ws.Range("D18") "3.142.98"
Set-TextValue $ws.Range("E18") "  -1.25%  "

# Row 19
Set-TextValue $ws.Range("D19") "3.41"
Set-TextValue $ws.Range("E19") "  +6.12%  "

# Row 20
Set-TextValue $ws.Range("D20") "0.0000226"
Set-TextValue $ws.Range("E20") "  +20.41%  "

# Row 21
Set-TextValue $ws.Range("D21") "13.23"
Set-TextValue $ws.Range("E21") "  -0.83%  "

# Row 22
Set-TextValue $ws.Range("D22") "424.75"
Set-TextValue $ws.Range("E22") "  -1.69%  "

# Row 23
Set-TextValue $ws.Range("D23") "8.44"
Set-TextValue $ws.Range("E23") "  -0.51%  "

# Row 24
Set-TextValue $ws.Range("D24") "4.91"
Set-TextValue $ws.Range("E24") "  -2.90%  "

# Row 25
Set-TextValue $ws.Range("E25") "  +6.42%  "

# Row 26
Set-TextValue $ws.Range("D26") "82.46"
Set-TextValue $ws.Range("E26") "  +10.75%  "

# Row 27
Set-TextValue $ws.Range("D27") "11.48"
Set-TextValue $ws.Range("E27") "  -2.35%  "

# Row 28
Set-TextValue $ws.Range("D28") "3.298.78"
Set-TextValue $ws.Range("E28") "  -3.89%  "

# Row 29
Set-TextValue $ws.Range("D29") "0.999"
Set-TextValue $ws.Range("E29") "  -0.10%  "

# Row 30
$ws.Range("B30").Value = "Binance-PegBSC-USD"
$ws.Range("C30").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextValue $ws.Range("D30") "0.999"
Set-TextValue $ws.Range("E30") "  -0.14%  "

# Row 31
$ws.Range("B31").Value = "Cronos"
$ws.Range("C31").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D31") "0.158"
Set-TextValue $ws.Range("E31") "  -5.90%  "

# Row 32
Set-TextValue $ws.Range("D32") "3.97"
Set-TextValue $ws.Range("E32") "  -2.55%  "

# Row 33
Set-TextValue $ws.Range("D33") "8.16"
Set-TextValue $ws.Range("E33") "  -3.29%  "

# Row 34
Set-TextValue $ws.Range("D34") "504.20"
Set-TextValue $ws.Range("E34") "  -5.28%  "

# Row 35
Set-TextValue $ws.Range("E35") "  +16.70%  "

# Row 36
Set-TextValue $ws.Range("D36") "6.94"
Set-TextValue $ws.Range("E36") "  +1.50%  "

# Row 37
Set-TextValue $ws.Range("E37") "  +3.41%  "

# Row 38
Set-TextValue $ws.Range("D38") "1.83"
Set-TextValue $ws.Range("E38") "  -1.41%  "

# Row 39
Set-TextValue $ws.Range("D39") "22.33"
Set-TextValue $ws.Range("E39") "  +2.20%  "

# Row 40
Set-TextValue $ws.Range("D40") "22.28"
Set-TextValue $ws.Range("E40") "  -0.02%  "

# Row 41
Set-TextValue $ws.Range("E41") "  +0.38%  "

# Row 42
Set-TextValue $ws.Range("E42") "  +0.08%  "

# Row 43
$ws.Range("B43").Value = "PolygonEcosystemToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
Set-TextValue $ws.Range("D43") "0.366"
Set-TextValue $ws.Range("E43") "  -2.03%  "

# Row 44
$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws.Range("D44") "1.86"
Set-TextValue $ws.Range("E44") "  -2.81%  "

# Row 45
$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue $ws.Range("D45") "0.134"
Set-TextValue $ws.Range("E45") "  +9.40%  "

# Row 46
$ws.Range("B46").Value = "Monero"
$ws.Range("C46").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D46") "145.74"
Set-TextValue $ws.Range("E46") "  -0.26%  "

# Row 47
Set-TextValue $ws.Range("D47") "43.71"
Set-TextValue $ws.Range("E47") "  -1.84%  "

# Row 48
Set-TextValue $ws.Range("D48") "164.02"
Set-TextValue $ws.Range("E48") "  -5.05%  "

# Row 49
Set-TextValue $ws.Range("D49") "0.0650"
Set-TextValue $ws.Range("E49") "  +11.80%  "

# Row 50
Set-TextValue $ws.Range("D50") "0.722"
Set-TextValue $ws.Range("E50") "  +2.81%  "

# Row 51
Set-TextValue $ws.Range("D51") "24.13"
Set-TextValue $ws.Range("E51") "  +0.67%  "
